$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - corrected values (serial parsing fix produced properly-parsed
# readings instead of the old placeholder/garbage values)
$ws.Range("A2").Value = "'0.14"
$ws.Range("B2").Value = "'-0.02"
$ws.Range("C2").Value = "'-0.00"
$ws.Range("D2").Value = "'16-Jan-2023 23:04:07"
$ws.Range("E2").Value = "'24.47"
$ws.Range("F2").Value = "'-0.02"
$ws.Range("G2").Value = "'-0.47"
$ws.Range("H2").Value = "'16-Jan-2023 23:04:07"

# Row 3 - new reading appended
$ws.Range("A3").Value = "'0.11"
$ws.Range("B3").Value = "'-0.03"
$ws.Range("C3").Value = "'-0.00"
$ws.Range("D3").Value = "'16-Jan-2023 23:05:08"
$ws.Range("E3").Value = "'24.47"
$ws.Range("F3").Value = "'-0.01"
$ws.Range("G3").Value = "'-0.43999999999999995"
$ws.Range("H3").Value = "'16-Jan-2023 23:05:08"

# Row 4 - new reading appended
$ws.Range("A4").Value = "'0.14"
$ws.Range("B4").Value = "'-0.02"
$ws.Range("C4").Value = "'-0.00"
$ws.Range("D4").Value = "'16-Jan-2023 23:06:08"
$ws.Range("E4").Value = "'24.47"
$ws.Range("F4").Value = "'-0.02"
$ws.Range("G4").Value = "'-0.54"
$ws.Range("H4").Value = "'16-Jan-2023 23:06:08"

# The leading apostrophes force Excel to keep these values as literal text
# (matching the source inlineStr cells) instead of auto-converting them to
# numbers/dates. Strip the resulting "quote prefix" formatting so the cell
# styles stay identical to the unstyled data rows in the original sheet.
$ws.Range("A2:H4").ClearFormats()
